$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-57 from 2023-09-20 (45189) to 2023-09-21 (45190)
for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 3).Value = 45190
}
